# Apply the "days" data correction described in the commit message:
# The values in column C (days) were originally measured as "days of life"
# since sowing. Since the two species were sown on different dates, the
# author re-based the values to the treatment period days by subtracting
# 12 from every value in column C for data rows 2 through 241.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 241; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = $cell.Value2 - 12
}
